# Insert a new weekly data row at row 117 (shifts old rows 117-175 down to 118-176)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("117:117").Insert()

$ws.Range("A117").Value = 9
$ws.Range("B117").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C117").Value = "Metropolitana"
$ws.Range("D117").Value = 44518
$ws.Range("D117").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E117").Value = 13
$ws.Range("F117").Value = 100112026
$ws.Range("G117").Value = "Haba"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 52
$ws.Range("K117").Value = 6000
$ws.Range("L117").Value = 7000
$ws.Range("M117").Value = 6500
$ws.Range("N117").Value = "`$/saco 25 kilos"
$ws.Range("O117").Value = "Región Metropolitana"
$ws.Range("P117").Value = 260
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"
